$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53 - this shifts the existing row 53
# (and everything below it) down by one, and grows the used range
# from A1:R177 to A1:R178 automatically.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new data point.
$ws.Cells.Item(53, 1).Value = 5
$ws.Cells.Item(53, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(53, 3).Value = "Maule"
$ws.Cells.Item(53, 4).Value = "2022-05-27"
$ws.Cells.Item(53, 5).Value = 7
$ws.Cells.Item(53, 6).Value = 100112017
$ws.Cells.Item(53, 7).Value = "Apio"
$ws.Cells.Item(53, 8).Value = "Americana (o)"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 700
$ws.Cells.Item(53, 11).Value = 6000
$ws.Cells.Item(53, 12).Value = 6000
$ws.Cells.Item(53, 13).Value = 6000
$ws.Cells.Item(53, 14).Value = "$/docena de matas"
$ws.Cells.Item(53, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(53, 16).Value = 1000
$ws.Cells.Item(53, 17).Value = 6
$ws.Cells.Item(53, 18).Value = "Hortaliza"
